$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($hf, $newName) {
    $count = $hf.Range.InlineShapes.Count
    if ($count -gt 0) {
        $s = $hf.Range.InlineShapes.Item(1)
        $s.ConvertToShape()
        $sh = $hf.Shapes.Item($hf.Shapes.Count)
        $sh.Name = $newName
        $sh.ConvertToInlineShape()
    }
}

# Default footer (appears on all pages except the first) -> Pearson logo: image1.png -> image2.png
$ftrDefault = $sec.Footers.Item(1)
Rename-InlinePicture $ftrDefault "image2.png"

# First-page footer -> Pearson logo: image1.png -> image2.png
$ftrFirst = $sec.Footers.Item(2)
Rename-InlinePicture $ftrFirst "image2.png"

# First-page header -> BTEC logo: image2.jpg -> image1.jpg
$hdrFirst = $sec.Headers.Item(2)
Rename-InlinePicture $hdrFirst "image1.jpg"
